# Insert a new data row at row 155 (pushing the existing rows 155-210 down
# to 156-211) and populate the new row with the "Fruta / hortaliza, semanal"
# record added in this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 155..210 down by one row, carrying formatting (incl. the date
# number format in column D) from the row above, as Excel does natively.
$ws.Rows(155).Insert()

# Populate the newly inserted row 155 with its values.
$ws.Cells.Item(155, 1).Value  = 1
$ws.Cells.Item(155, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(155, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(155, 4).Value  = 44636
$ws.Cells.Item(155, 5).Value  = 15
$ws.Cells.Item(155, 6).Value  = "Fruta"
$ws.Cells.Item(155, 7).Value  = 100102
$ws.Cells.Item(155, 8).Value  = "Cítricos"
$ws.Cells.Item(155, 9).Value  = 100102003
$ws.Cells.Item(155, 10).Value = "Limón"
$ws.Cells.Item(155, 11).Value = "Sin especificar"
$ws.Cells.Item(155, 12).Value = "2a amarillo"
$ws.Cells.Item(155, 13).Value = 300
$ws.Cells.Item(155, 14).Value = 32000
$ws.Cells.Item(155, 15).Value = 35000
$ws.Cells.Item(155, 16).Value = 33500
$ws.Cells.Item(155, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(155, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(155, 19).Value = 1675
$ws.Cells.Item(155, 20).Value = 20
